$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so values like "3.70" or "98.033.78"
# are not auto-coerced into numbers (which would drop trailing zeros / misparse multi-dot values).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.033.78"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.430.38"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "257.95"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "659.78"
$ws.Range("E6").Value = "  +5.68%  "
$ws.Range("D7").Value = "1.50"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +6.28%  "
$ws.Range("D9").Value = "1.06"
$ws.Range("E9").Value = "  +10.09%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.425.86"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("E12").Value = "  +6.70%  "
$ws.Range("D13").Value = "42.23"
$ws.Range("E13").Value = "  +6.22%  "
$ws.Range("D14").Value = "6.50"
$ws.Range("E14").Value = "  +18.83%  "
$ws.Range("E15").Value = "  +4.19%  "
$ws.Range("D16").Value = "97.819.93"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "4.070.91"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "8.66"
$ws.Range("E18").Value = "  +36.91%  "
$ws.Range("D19").Value = "3.406.33"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  +15.00%  "
$ws.Range("E21").Value = "  +65.35%  "
$ws.Range("D22").Value = "11.10"
$ws.Range("E22").Value = "  +18.32%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "511.24"
$ws.Range("E24").Value = "  +4.81%  "
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("D26").Value = "6.19"
$ws.Range("E26").Value = "  +9.86%  "
$ws.Range("D27").Value = "99.25"
$ws.Range("E27").Value = "  +11.45%  "
$ws.Range("D28").Value = "12.89"
$ws.Range("E28").Value = "  +7.72%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.609.45"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.155"
$ws.Range("E30").Value = "  +13.41%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "11.49"
$ws.Range("E31").Value = "  +11.17%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").Value = "0.198"
$ws.Range("E32").Value = "  +5.19%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "0.576"
$ws.Range("E35").Value = "  +22.33%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "29.97"
$ws.Range("E36").Value = "  +7.62%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "2.24"
$ws.Range("E37").Value = "  +15.20%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "7.88"
$ws.Range("E38").Value = "  +9.59%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.157"
$ws.Range("E39").Value = "  +5.98%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.42"
$ws.Range("E40").Value = "  +15.33%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "519.17"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "24.74"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.861"
$ws.Range("E43").Value = "  +9.59%  "
$ws.Range("D44").Value = "0.0423"
$ws.Range("E44").Value = "  +27.10%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "3.70"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "5.47"
$ws.Range("E47").Value = "  +15.45%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "8.22"
$ws.Range("E48").Value = "  +12.49%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "1.60"
$ws.Range("E50").Value = "  +18.00%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  +7.85%  "

Write-Host "Update complete"
